$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new advice texts for "Thiên Mã" section (rows 13-15, column B)
# Set in this order so new shared strings are appended matching the source order
$ws.Range("B15").Value = "Bạn đi ra ngoài hay gặp điều kém may mắn nhưng không vì thế mà không cố gắng, chỉ có năng động tích cực mới giúp bạn có nhiều cơ hội hơn."
$ws.Range("B14").Value = "Bạn ra ngoài luôn gặp may mắn, tiến hành công việc dễ xứng ý toại lòng. Nhưng không vì thế mà bất cẩn trong đi lại hay công việc."
$ws.Range("B13").Value = "Bạn có thể gặp khó khăn khi đi lại hoặc trong công việc nhưng hầu hết bạn đều có thể xử lý tốt mọi vấn đề phát sinh."

# Update the active cell selection to match the saved view state
$ws.Range("H28").Select()
